$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$tbl = $ws.ListObjects.Item(1)

# 1. Pre-expand the table by one row so the bottom data row gets copied/
#    pushed into the new last row (row 196) with its calculated-column
#    formula intact when we shift everything down.
$tbl.Resize($ws.Range("A8:K196"))

# 2. Insert a brand-new blank row at row 96 (CTO entry), pushing every
#    row from 96 downward (old row 96 -> 97, ..., old row 195 -> 196).
$ws.Rows("96:96").Insert()

# 3. Inherit formatting for the new row 96 from the row above it (row 95)
#    so its style indices match the surrounding table rows.
$ws.Range("A95:K95").Copy()
$ws.Range("A96:K96").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 4. The new row's calculated "EARNED " column lost its formula during the
#    insert (a blank inserted row isn't auto-filled by the table engine),
#    so restore it explicitly to match the other rows.
$ws.Range("G96").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

# 5. Give the three new leave-card dates (K94:K96) the same date format
#    already used by the other REMARKS date cells in this column.
$ws.Range("K90").Copy()
$ws.Range("K94:K96").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 6. Fill in the new CTO leave entries.
$ws.Range("B94").Value = "SL(1-0-0)"
$ws.Range("C94").Value = 1.25
$ws.Range("H94").Value = 1
$ws.Range("K94").Value = "10/5/2023"

$ws.Range("B95").Value = "FL(1-0-0)"
$ws.Range("D95").Value = 1
$ws.Range("K95").Value = "11/3/2023"

$ws.Range("B96").Value = "VL(1-0-0)"
$ws.Range("D96").Value = 1
$ws.Range("K96").Value = "11/17/2023"

# 7. Match the author's last on-screen selection.
$ws.Range("K96").Select()
